$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cell K1 "Gender" - copy formatting from the adjacent header cell J1
# so it picks up the same bold/centered/bordered header style, then set its text.
$ws.Range("J1").Copy($ws.Range("K1"))
$ws.Range("K1").Value = "Gender"

# Gender formula column, mirroring the existing GenPos/CatPos formula columns.
# Row 2 (top result, outside the two shared-formula blocks used by column J).
$ws.Range("K2").Formula = '=IF(LEFT(D2,1)="M","M","F")'

# Rows 3:66 share one formula (first male/female results block).
$ws.Range("K3:K66").Formula = '=IF(LEFT(D3,1)="M","M","F")'

# Rows 67:86 share a second formula block (mirrors J67:J86's shared formula split).
$ws.Range("K67:K86").Formula = '=IF(LEFT(D67,1)="M","M","F")'

# Match the saved selection state recorded in the edited workbook.
[void]$ws.Range("K2").Select()
